# Update the "Corr/total marks" on the marksheet:
#  - Marking row (B11): 3 -> 5
#  - Total row (B12): 84 -> 140, with the matching "Max" display text (E12) updated too

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 140
$ws.Range("E12").Value = "140/140"
